$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://docs.github.com/es/get-started/writing-on-github/getting-started-with-writing-and-formatting-on-github/basic-writing-and-formatting-syntax"
$desc = "Documentación Sintaxis MK archivo readme github"

# Add the new row of data (URL + description) mirroring the existing rows
$ws.Range("B20").Value = $url
$ws.Range("C20").Value = $desc

# Turn B20 into a real hyperlink (adds a relationship + hyperlinks entry)
$ws.Hyperlinks.Add($ws.Range("B20"), $url)

# Restore the "Hipervínculo" look (underlined themed link style) used by the
# other link cells in column B, since Hyperlinks.Add applies its own style
$ws.Range("B20").Style = "Hipervínculo"

# Move the active selection down to reflect the new last row, like Excel
# would after typing into row 20 and moving to the next row
$ws.Range("C21").Select()
